$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (sharedStrings via cell writes) ---
# "Volume 31   Number  41" -> "...Number  42" (A8 holds the rich-text run)
$ws.Range("A8").Value = "Volume 31   Number  42"
# "Report Covering the Week  10/7/2024  Through  10/13/2024" -> 10/14/2024 .. 10/20/2024 (C9)
$ws.Range("C9").Value = "Report Covering the Week  10/14/2024  Through  10/20/2024"

# --- Data table updates (rows 15-30) ---

# Row 15
$ws.Range("A14").Copy() | Out-Null
$ws.Range("F15").PasteSpecial(-4122) | Out-Null
$ws.Range("F15").Value = "0"
$ws.Range("A14").Copy() | Out-Null
$ws.Range("G15").PasteSpecial(-4122) | Out-Null
$ws.Range("G15").Value = "0"
$ws.Range("A14").Copy() | Out-Null
$ws.Range("H15").PasteSpecial(-4122) | Out-Null
$ws.Range("H15").Value = "***.*"
$ws.Range("N15").Value = -50.0

# Row 16
$ws.Range("A14").Copy() | Out-Null
$ws.Range("D16").PasteSpecial(-4122) | Out-Null
$ws.Range("D16").Value = "0"
$ws.Range("A14").Copy() | Out-Null
$ws.Range("E16").PasteSpecial(-4122) | Out-Null
$ws.Range("E16").Value = "***.*"
$ws.Range("F16").Value = 3.0
$ws.Range("G16").Value = 7.0
$ws.Range("H16").Value = -57.142857142857
$ws.Range("I16").Value = 51.0
$ws.Range("K16").Value = -7.272727272727
$ws.Range("L16").Value = 4.081632653061
$ws.Range("M16").Value = -40.0
$ws.Range("N16").Value = -83.441558441558

# Row 17
$ws.Range("C17").Value = 3.0
$ws.Range("D17").Value = 1.0
$ws.Range("E17").Value = 200.0
$ws.Range("F17").Value = 9.0
$ws.Range("G17").Value = 3.0
$ws.Range("H17").Value = 200.0
$ws.Range("I17").Value = 100.0
$ws.Range("J17").Value = 87.0
$ws.Range("K17").Value = 14.942528735632
$ws.Range("L17").Value = 26.582278481012
$ws.Range("M17").Value = 28.205128205128
$ws.Range("N17").Value = -60.629921259842

# Row 18
$ws.Range("F18").Value = 3.0
$ws.Range("G18").Value = 2.0
$ws.Range("H18").Value = 50.0
$ws.Range("J18").Value = 76.0
$ws.Range("K18").Value = -5.263157894736
$ws.Range("L18").Value = -18.181818181818
$ws.Range("M18").Value = -19.101123595505
$ws.Range("N18").Value = -78.823529411764

# Row 19
$ws.Range("D19").Value = 4.0
$ws.Range("E19").Value = -50.0
$ws.Range("F19").Value = 9.0
$ws.Range("G19").Value = 13.0
$ws.Range("H19").Value = -30.76923076923
$ws.Range("I19").Value = 165.0
$ws.Range("J19").Value = 132.0
$ws.Range("K19").Value = 25.0
$ws.Range("L19").Value = 0.0
$ws.Range("M19").Value = -11.290322580645
$ws.Range("N19").Value = -15.384615384615

# Row 20
$ws.Range("A14").Copy() | Out-Null
$ws.Range("C20").PasteSpecial(-4122) | Out-Null
$ws.Range("C20").Value = "0"
$ws.Range("E20").Value = -100.0
$ws.Range("F20").Value = 7.0
$ws.Range("H20").Value = 75.0
$ws.Range("J20").Value = 51.0
$ws.Range("K20").Value = -1.960784313725
$ws.Range("L20").Value = 6.382978723404
$ws.Range("M20").Value = -3.846153846153
$ws.Range("N20").Value = -84.709480122324

# Row 21
$ws.Range("C21").Value = 6.0
$ws.Range("D21").Value = 7.0
$ws.Range("E21").Value = -14.285714285714
$ws.Range("F21").Value = 31.0
$ws.Range("G21").Value = 29.0
$ws.Range("H21").Value = 6.896551724137
$ws.Range("I21").Value = 443.0
$ws.Range("J21").Value = 407.0
$ws.Range("K21").Value = 8.845208845208
$ws.Range("L21").Value = 2.073732718894
$ws.Range("M21").Value = -10.141987829614
$ws.Range("N21").Value = -69.236111111111

# Row 22
$ws.Range("M22").Value = -83.333333333333

# Row 23
$ws.Range("C23").Value = 3.0
$ws.Range("A14").Copy() | Out-Null
$ws.Range("D23").PasteSpecial(-4122) | Out-Null
$ws.Range("D23").Value = "0"
$ws.Range("A14").Copy() | Out-Null
$ws.Range("E23").PasteSpecial(-4122) | Out-Null
$ws.Range("E23").Value = "***.*"
$ws.Range("F23").Value = 7.0
$ws.Range("G23").Value = 3.0
$ws.Range("H23").Value = 133.333333333333
$ws.Range("I23").Value = 88.0
$ws.Range("K23").Value = -1.123595505617
$ws.Range("L23").Value = -7.368421052631
$ws.Range("M23").Value = 29.411764705882

# Row 24
$ws.Range("C24").Value = 18.0
$ws.Range("E24").Value = 80.0
$ws.Range("F24").Value = 51.0
$ws.Range("H24").Value = -10.526315789473
$ws.Range("I24").Value = 432.0
$ws.Range("J24").Value = 494.0
$ws.Range("K24").Value = -12.550607287449
$ws.Range("L24").Value = 16.44204851752
$ws.Range("M24").Value = 18.356164383561

# Row 25
$ws.Range("C25").Value = 4.0
$ws.Range("D25").Value = 5.0
$ws.Range("E25").Value = -20.0
$ws.Range("F25").Value = 23.0
$ws.Range("G25").Value = 36.0
$ws.Range("H25").Value = -36.111111111111
$ws.Range("I25").Value = 215.0
$ws.Range("J25").Value = 255.0
$ws.Range("K25").Value = -15.686274509803
$ws.Range("L25").Value = 119.387755102041

# Row 26
$ws.Range("C26").Value = 5.0
$ws.Range("A14").Copy() | Out-Null
$ws.Range("D26").PasteSpecial(-4122) | Out-Null
$ws.Range("D26").Value = "0"
$ws.Range("A14").Copy() | Out-Null
$ws.Range("E26").PasteSpecial(-4122) | Out-Null
$ws.Range("E26").Value = "***.*"
$ws.Range("G26").Value = 7.0
$ws.Range("H26").Value = 185.714285714286
$ws.Range("I26").Value = 132.0
$ws.Range("K26").Value = 4.761904761904
$ws.Range("L26").Value = -10.204081632653
$ws.Range("M26").Value = -43.103448275862

# Row 27
$ws.Range("A14").Copy() | Out-Null
$ws.Range("F27").PasteSpecial(-4122) | Out-Null
$ws.Range("F27").Value = "0"
$ws.Range("G27").Value = 1.0
$ws.Range("H27").Value = -100.0

# Row 28
$ws.Range("J14").Copy() | Out-Null
$ws.Range("D28").PasteSpecial(-4122) | Out-Null
$ws.Range("D28").Value = 1.0
$ws.Range("K14").Copy() | Out-Null
$ws.Range("E28").PasteSpecial(-4122) | Out-Null
$ws.Range("E28").Value = -100.0
$ws.Range("J14").Copy() | Out-Null
$ws.Range("G28").PasteSpecial(-4122) | Out-Null
$ws.Range("G28").Value = 1.0
$ws.Range("K14").Copy() | Out-Null
$ws.Range("H28").PasteSpecial(-4122) | Out-Null
$ws.Range("H28").Value = -100.0
$ws.Range("J28").Value = 11.0
$ws.Range("K28").Value = 9.090909090909

# Row 29
$ws.Range("A14").Copy() | Out-Null
$ws.Range("G29").PasteSpecial(-4122) | Out-Null
$ws.Range("G29").Value = "0"
$ws.Range("A14").Copy() | Out-Null
$ws.Range("H29").PasteSpecial(-4122) | Out-Null
$ws.Range("H29").Value = "***.*"

# Row 30
$ws.Range("A14").Copy() | Out-Null
$ws.Range("G30").PasteSpecial(-4122) | Out-Null
$ws.Range("G30").Value = "0"
$ws.Range("A14").Copy() | Out-Null
$ws.Range("H30").PasteSpecial(-4122) | Out-Null
$ws.Range("H30").Value = "***.*"
